$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Rows.Item(11).Delete()
